$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record (weekly update) before the current row 113 — this
# shifts the existing rows 113:136 down to 114:137 and grows the used
# range to A1:R137.
$ws.Rows.Item(113).Insert()

# Populate the newly inserted row 113 with the new "Sandia" price record.
$ws.Cells.Item(113, 1).Value = 4
$ws.Cells.Item(113, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(113, 3).Value = "Los Lagos"
$ws.Cells.Item(113, 4).Value = 44511
$ws.Cells.Item(113, 5).Value = 10
$ws.Cells.Item(113, 6).Value = 100112028
$ws.Cells.Item(113, 7).Value = "Sandia"
$ws.Cells.Item(113, 8).Value = "Sin especificar"
$ws.Cells.Item(113, 9).Value = "Primera"
$ws.Cells.Item(113, 10).Value = 300
$ws.Cells.Item(113, 11).Value = 1200
$ws.Cells.Item(113, 12).Value = 1200
$ws.Cells.Item(113, 13).Value = 1200
$ws.Cells.Item(113, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(113, 15).Value = "Perú"
$ws.Cells.Item(113, 16).Value = 1200
$ws.Cells.Item(113, 17).Value = 1
$ws.Cells.Item(113, 18).Value = "Hortaliza"
